$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.222.05"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.82"
$ws.Range("E3").Value = "  +4.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.92"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4439"
$ws.Range("E7").Value = "  +5.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3698"
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07693"
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.126"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.07"
$ws.Range("E13").Value = "  +3.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.267"
$ws.Range("E14").Value = "  +3.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.573"
$ws.Range("E15").Value = "  +5.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.841.03"
$ws.Range("E16").Value = "  +6.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.47"
$ws.Range("E17").Value = "  +6.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001084"
$ws.Range("E18").Value = "  +1.98%  "
$ws.Range("E19").Value = "  +9.75%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.51"
$ws.Range("E21").Value = "  +4.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.189"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.282.27"
$ws.Range("E23").Value = "  +3.45%  "
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.083"
$ws.Range("E25").Value = "  -13.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.78"
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.88"
$ws.Range("E27").Value = "  +3.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.036.87"
$ws.Range("E28").Value = "  +5.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.321"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.33"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.208"
$ws.Range("E31").Value = "  -5.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.868"
$ws.Range("E32").Value = "  +5.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09214"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.665"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02353"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2179"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06213"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6569"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.150"
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.197"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.163"
$ws.Range("E42").Value = "  +3.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6074"
$ws.Range("E46").Value = "  +4.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.766"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.49"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.038"
$ws.Range("E49").Value = "  +5.33%  "
$ws.Range("E50").Value = "  +5.37%  "
$ws.Range("E51").Value = "  +2.09%  "
